$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value2 = 1318.3334
$ws.Range("I34").Value2 = 1318.3334
$ws.Range("K34").Value2 = 1318.3334
$ws.Range("M34").Value2 = -1115.3334

$ws.Range("H36").Value2 = 1318.3334
$ws.Range("I36").Value2 = 1318.3334
$ws.Range("K36").Value2 = 1318.3334
$ws.Range("M36").Value2 = -603.3334

$ws.Range("H41").Value2 = 425
$ws.Range("I41").Value2 = 425
$ws.Range("K41").Value2 = 425
$ws.Range("M41").Value2 = 15

$ws.Range("H53").Value2 = 555.8125
$ws.Range("I53").Value2 = 417.63635
$ws.Range("K53").Value2 = 417.63635
$ws.Range("M53").Value2 = 219.36365

$ws.Range("H62").Value2 = 11500
$ws.Range("I62").Value2 = 9500
$ws.Range("K62").Value2 = 9500
$ws.Range("M62").Value2 = -8876

$ws.Range("H65").Value2 = 11500
$ws.Range("I65").Value2 = 9500
$ws.Range("K65").Value2 = 47500
$ws.Range("M65").Value2 = -44380

$ws.Range("H86").Value2 = 3000
$ws.Range("I86").Value2 = 3000
$ws.Range("J86").Value2 = 0
$ws.Range("K86").Value2 = 3000
$ws.Range("L86").Value2 = 0
$ws.Range("M86").Value2 = -1877
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value2 = 3000
$ws.Range("I89").Value2 = 3000
$ws.Range("J89").Value2 = 0
$ws.Range("K89").Value2 = 15000
$ws.Range("L89").Value2 = 0
$ws.Range("M89").Value2 = -9384
$ws.Range("N89").ClearContents()

$ws.Range("H96").Value2 = 598.8
$ws.Range("I96").Value2 = 258.33334
$ws.Range("J96").Value2 = 1109.5
$ws.Range("K96").Value2 = 775.0000200000001
$ws.Range("L96").Value2 = 3328.5
$ws.Range("M96").Value2 = 597.9999799999999
$ws.Range("N96").Value2 = -6074.5

$ws.Range("H106").Value2 = 5000
$ws.Range("I106").Value2 = 5000
$ws.Range("K106").Value2 = 5000
$ws.Range("M106").Value2 = -4369

$ws.Range("H113").Value2 = 4122.5
$ws.Range("I113").Value2 = 4122.5
$ws.Range("J113").Value2 = 0
$ws.Range("K113").Value2 = 4122.5
$ws.Range("L113").Value2 = 0
$ws.Range("M113").Value2 = -868.5
$ws.Range("N113").ClearContents()

$ws.Range("H116").Value2 = 17000
$ws.Range("I116").Value2 = 17000
$ws.Range("K116").Value2 = 17000
$ws.Range("M116").Value2 = -13558

$ws.Range("H132").Value2 = 7364.8887
$ws.Range("I132").Value2 = 7364.8887
$ws.Range("K132").Value2 = 22094.6661
$ws.Range("M132").Value2 = -19564.6661

$ws.Range("H138").Value2 = 2999
$ws.Range("J138").Value2 = 4542.857
$ws.Range("L138").Value2 = 13628.571
$ws.Range("N138").Value2 = -23908.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H95").Value2 = 31772.715
$ws.Range("J95").Value2 = 31772.715
$ws.Range("L95").Value2 = 31772.715
$ws.Range("N95").Value2 = -37264.715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H139").Value2 = 100000
$ws.Range("J139").Value2 = 100000
$ws.Range("L139").Value2 = 100000
$ws.Range("N139").Value2 = -110280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value2 = 0
$ws.Range("I132").Value2 = 0
$ws.Range("K132").Value2 = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value2 = 29.666666
$ws.Range("I6").Value2 = 29.666666
$ws.Range("J6").Value2 = 0
$ws.Range("K6").Value2 = 88.99999800000001
$ws.Range("L6").Value2 = 0
$ws.Range("M6").Value2 = 24.00000199999999
$ws.Range("N6").ClearContents()

$ws.Range("H131").Value2 = 1800
$ws.Range("J131").Value2 = 1800
$ws.Range("L131").Value2 = 5400
$ws.Range("N131").Value2 = -15480

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value2 = 67500
$ws.Range("J39").Value2 = 67500
$ws.Range("L39").Value2 = 67500
$ws.Range("N39").Value2 = -68564

$ws.Range("H98").Value2 = 70993.336
$ws.Range("J98").Value2 = 70993.336
$ws.Range("L98").Value2 = 70993.336
$ws.Range("N98").Value2 = -76983.336

$ws.Range("H132").Value2 = 0
$ws.Range("I132").Value2 = 0
$ws.Range("K132").Value2 = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 1969.7059
$ws.Range("I7").Value2 = 1717.8125
$ws.Range("K7").Value2 = 1717.8125
$ws.Range("M7").Value2 = -1605.8125

$ws.Range("H40").Value2 = 3002
$ws.Range("I40").Value2 = 3002
$ws.Range("K40").Value2 = 3002
$ws.Range("M40").Value2 = -2866

$ws.Range("H46").Value2 = 625
$ws.Range("J46").Value2 = 666.6667
$ws.Range("L46").Value2 = 666.6667
$ws.Range("N46").Value2 = -1042.6667

$ws.Range("H55").Value2 = 543
$ws.Range("I55").Value2 = 543
$ws.Range("J55").Value2 = 0
$ws.Range("K55").Value2 = 543
$ws.Range("L55").Value2 = 0
$ws.Range("M55").Value2 = -370
$ws.Range("N55").ClearContents()

$ws.Range("H93").Value2 = 0
$ws.Range("I93").Value2 = 0
$ws.Range("K93").Value2 = 0
$ws.Range("M93").ClearContents()

$ws.Range("H126").Value2 = 1969.7059
$ws.Range("I126").Value2 = 1717.8125
$ws.Range("K126").Value2 = 5153.4375
$ws.Range("M126").Value2 = -2683.4375

$ws.Range("H132").Value2 = 11999.5
$ws.Range("I132").Value2 = 11999.5
$ws.Range("K132").Value2 = 35998.5
$ws.Range("M132").Value2 = -33468.5

$ws.Range("H136").Value2 = 1800
$ws.Range("I136").Value2 = 1800
$ws.Range("J136").Value2 = 0
$ws.Range("K136").Value2 = 5400
$ws.Range("L136").Value2 = 0
$ws.Range("M136").Value2 = -2850
$ws.Range("N136").ClearContents()

$ws.Range("H139").Value2 = 45789
$ws.Range("I139").Value2 = 45789
$ws.Range("K139").Value2 = 45789
$ws.Range("M139").Value2 = -40649

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value2 = 60114
$ws.Range("I64").Value2 = 0
$ws.Range("J64").Value2 = 60114
$ws.Range("K64").Value2 = 0
$ws.Range("L64").Value2 = 60114
$ws.Range("N64").Value2 = -60610
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value2 = 60114
$ws.Range("I67").Value2 = 0
$ws.Range("J67").Value2 = 60114
$ws.Range("K67").Value2 = 0
$ws.Range("L67").Value2 = 60114
$ws.Range("N67").Value2 = -61830
$ws.Range("M67").ClearContents()

$ws.Range("H76").Value2 = 40000
$ws.Range("J76").Value2 = 40000
$ws.Range("L76").Value2 = 40000
$ws.Range("N76").Value2 = -40630

$ws.Range("H79").Value2 = 40000
$ws.Range("J79").Value2 = 40000
$ws.Range("L79").Value2 = 40000
$ws.Range("N79").Value2 = -42184

$ws.Range("H81").Value2 = 879.8333
$ws.Range("I81").Value2 = 879.8333
$ws.Range("J81").Value2 = 0
$ws.Range("K81").Value2 = 1759.6666
$ws.Range("L81").Value2 = 0
$ws.Range("M81").Value2 = -698.6666
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value2 = 879.8333
$ws.Range("I84").Value2 = 879.8333
$ws.Range("J84").Value2 = 0
$ws.Range("K84").Value2 = 8798.333000000001
$ws.Range("L84").Value2 = 0
$ws.Range("M84").Value2 = -3494.333000000001
$ws.Range("N84").ClearContents()

$ws.Range("H123").Value2 = 275000
$ws.Range("J123").Value2 = 275000
$ws.Range("L123").Value2 = 275000
$ws.Range("N123").Value2 = -284800

$ws.Range("H126").Value2 = 3046
$ws.Range("I126").Value2 = 2766.8572
$ws.Range("K126").Value2 = 8300.571599999999
$ws.Range("M126").Value2 = -5830.571599999999
